$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("samples_retained")

# Update the EmoReact_V_1.0 row (row 12) counts
$ws.Range("C12").Value = 579
$ws.Range("D12").Value = 295
$ws.Range("E12").Value = 38

# Update the note for that row, appending a pointer to Notes.docx
$ws.Range("H12").Value = "N. A. children's English; see Notes.docx"

# Update the selected cell to reflect the author's last position
$null = $ws.Range("C13").Select()
